# Update the stack-trace line numbers embedded in the document text to
# reflect the new M2DocEvaluator/M2DocUtils/AbstractTemplatesTestSuite
# source line numbers (caused by adding the M2Doc version to template
# custom properties, which shifted subsequent code).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "M2DocEvaluator.java:540)" "M2DocEvaluator.java:543)"
Replace-Text "M2DocEvaluator.java:1038)" "M2DocEvaluator.java:1084)"
Replace-Text "M2DocEvaluator.java:1254)" "M2DocEvaluator.java:1300)"
Replace-Text "M2DocEvaluator.java:275)" "M2DocEvaluator.java:278)"
Replace-Text "M2DocEvaluator.java:264)" "M2DocEvaluator.java:267)"
Replace-Text "M2DocUtils.java:712)" "M2DocUtils.java:694)"
Replace-Text "AbstractTemplatesTestSuite.java:459)" "AbstractTemplatesTestSuite.java:475)"
Replace-Text "AbstractTemplatesTestSuite.java:369)" "AbstractTemplatesTestSuite.java:384)"
